# Apply the "working for R joints! yay." edit to the Inputs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")
$ws.Activate()

# Row 3 (Joint #1): X coord unchanged request -> C3 0 -> 30
$ws.Range("C3").Value = 30

# Row 4 (Joint #2): C4 0 -> 30
$ws.Range("C4").Value = 30

# Row 5 (Joint #3): D5 20 -> 10, F5 "coupler,output" -> "coupler,intermed"
$ws.Range("D5").Value = 10
$ws.Range("F5").Value = "coupler,intermed"

# Row 6 (Joint #4): F6 "output 0" -> "intermed, gnd"
$ws.Range("F6").Value = "intermed, gnd"

# Row 7 (Joint #5): C7 15 -> 20, D7 25 -> 20, F7 "coupler" -> "intermed, coupler2"
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 20
$ws.Range("F7").Value = "intermed, coupler2"

# Row 8 (Joint #6): was fully blank, now populated with an R joint
$ws.Range("B8").Value = "r"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 20
$ws.Range("F8").Value = "coupler,output"

# Row 9 (Joint #7): was fully blank, now populated with an R joint
$ws.Range("B9").Value = "r"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = "output, gnd"

# Update the active selection to match the post-edit state (F10).
$ws.Range("F10").Select()

